$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; this pushes all existing data down by
# one row (A1 -> A2, A2 -> A3, ..., A2023 -> A2024) and creates a fresh,
# empty A1 for the new CSV header.
$ws.Rows("1:1").Insert()

# Write the header line into the newly created first row.
$ws.Range("A1").Value = "time,position,angle,motor_velocity,pendulum_velocity,voltage,energy,mode"

# Restore the selection to match the post-edit state.
$ws.Range("F15").Select()
